$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the added "canonical SMILES" column (D)
$ws.Range("D2").Value = "canonical SMILES"

# Give column D a width close to the target (42.28515625 chars of stored
# width); ColumnWidth only lands on a 1/7-character grid, so 41.5 is the
# closest input that converges on the nearest achievable stored width.
$ws.Columns.Item(4).ColumnWidth = 41.5

# For every data row, the new "canonical SMILES" value is simply the
# existing "canonical isomeric SMILES" value (column C) with the
# cis/trans stereo-bond slash markers ("/" and "\") stripped out.
for ($r = 3; $r -le 17; $r++) {
    $isomeric = $ws.Cells.Item($r, 3).Value2
    $canonical = $isomeric.Replace("/", "").Replace("\", "")
    $ws.Cells.Item($r, 4).Value = $canonical
}
